$d = $word.ActiveDocument
$d.Content.Find.Execute("Tresidder", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Tressider", 2)
